# QualityCheckLogBook.xlsx update
# Commit: "change link to quanlity check log book"
#
# Fills in the head-movement QC columns (H/I/J) for subjects 15-30 that were
# previously left blank, normalises the "1" placeholder answers in column J
# to the descriptive "none" answer used elsewhere in the sheet, expands the
# final open question about nuisance regressors, and hides the now-unused
# row 33 (sub031). Also re-establishes the freeze-panes view that was in
# place when the sheet was last looked at (header rows + subid column
# frozen, scrolled down near the bottom of the table).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Column J: replace placeholder "1" values with the descriptive "none" ---
$noneRows = @(3, 8, 9, 10, 11, 12, 13, 15)
foreach ($r in $noneRows) {
    $ws.Cells.Item($r, 10).Value = "none"
}

# --- Rows 17-32: fill in the previously-empty H/I/J QC columns ---
$ws.Cells.Item(17, 8).Value = 1
$ws.Cells.Item(17, 9).Value = 0.2
$ws.Cells.Item(17, 10).Value = "1 small shift in main task run1"

$ws.Cells.Item(18, 8).Value = 1
$ws.Cells.Item(18, 9).Value = 0.2
$ws.Cells.Item(18, 10).Value = "none"

$ws.Cells.Item(19, 8).Value = 1
$ws.Cells.Item(19, 9).Value = 0.2
$ws.Cells.Item(19, 10).Value = "none"

$ws.Cells.Item(20, 8).Value = 1
$ws.Cells.Item(20, 9).Value = 0.2
$ws.Cells.Item(20, 10).Value = "none"

$ws.Cells.Item(21, 8).Value = 1
$ws.Cells.Item(21, 9).Value = 0.2
$ws.Cells.Item(21, 10).Value = "none"

$ws.Cells.Item(22, 8).Value = 0.6
$ws.Cells.Item(22, 9).Value = 0.2
$ws.Cells.Item(22, 10).Value = "none"

$ws.Cells.Item(23, 8).Value = 1
$ws.Cells.Item(23, 9).Value = 0.2
$ws.Cells.Item(23, 10).Value = "none"

$ws.Cells.Item(24, 8).Value = 0.5
$ws.Cells.Item(24, 9).Value = 0.4
$ws.Cells.Item(24, 10).Value = "none"

$ws.Cells.Item(25, 8).Value = 2
$ws.Cells.Item(25, 9).Value = 0.2
$ws.Cells.Item(25, 10).Value = "2 in localizer task"

$ws.Cells.Item(26, 8).Value = 1
$ws.Cells.Item(26, 9).Value = 0.5
$ws.Cells.Item(26, 10).Value = "none"

$ws.Cells.Item(27, 8).Value = 0.5
$ws.Cells.Item(27, 9).Value = 0.2
$ws.Cells.Item(27, 10).Value = "none"

$ws.Cells.Item(28, 8).Value = 2
$ws.Cells.Item(28, 9).Value = 1
$ws.Cells.Item(28, 10).Value = "some in all runs"

$ws.Cells.Item(29, 8).Value = 1
$ws.Cells.Item(29, 9).Value = 0.2
$ws.Cells.Item(29, 10).Value = "none"

$ws.Cells.Item(30, 8).Value = 2
$ws.Cells.Item(30, 9).Value = 0.2
$ws.Cells.Item(30, 10).Value = "none"

$ws.Cells.Item(31, 8).Value = 1
$ws.Cells.Item(31, 9).Value = 0.2
$ws.Cells.Item(31, 10).Value = "none"

$ws.Cells.Item(32, 8).Value = 0.2
$ws.Cells.Item(32, 9).Value = 0.2
$ws.Cells.Item(32, 10).Value = "none"

# --- Row 33 (sub031) is no longer relevant - hide it ---
$ws.Rows.Item(33).Hidden = $true

# --- Expand the final open question about nuisance regressors ---
$ws.Range("H37").Value = "constructing nuisance regressors for sudden head movement: what count as sudden head movement (framewise displacement > voxel size?)"

# --- Restore the sheet's frozen-header view (subid column + the two header
#     rows), scrolled to the bottom of the table, with H38 as the active cell ---
$ws.Range("B3").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("H38").Select()
